$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (changed) date column C for rows 2-9 from 2023-10-25 (45224) to 2023-11-03 (45233)
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
